$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly export gained a new week's worth of records (4 rows) that must
# be inserted right after the header block of existing "Brocoli" data, at
# row 244 — pushing every subsequent data row down by 4 (old row N -> N+4).
$ws.Rows("244:247").Insert()

# New week's data (4 quality/origin combinations), inserted at rows 244-247.
$newRows = @(
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44455, 13, 100112023, "Brócoli", "Sin especificar", "Primera", 4300, 600, 650, 625, "$/unidad", "Región Metropolitana", 625, 1, "Hortaliza"),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44455, 13, 100112023, "Brócoli", "Sin especificar", "Primera", 3400, 600, 650, 625, "$/unidad", "Región de O'Higgins", 625, 1, "Hortaliza"),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44455, 13, 100112023, "Brócoli", "Sin especificar", "Segunda", 1600, 450, 500, 475, "$/unidad", "Región Metropolitana", 475, 1, "Hortaliza"),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44455, 13, 100112023, "Brócoli", "Sin especificar", "Segunda", 970, 450, 500, 475, "$/unidad", "Región de O'Higgins", 475, 1, "Hortaliza")
)

$startRow = 244
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
